$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 10).Value = 1.06
$ws.Cells.Item(2, 11).Value = 10
$ws.Cells.Item(2, 14).Value = 2.08
$ws.Cells.Item(2, 15).Value = 1.73
$ws.Cells.Item(3, 8).Value = 3
$ws.Cells.Item(3, 9).Value = 2.9
$ws.Cells.Item(3, 10).Value = 1.11
$ws.Cells.Item(3, 11).Value = 6.5
$ws.Cells.Item(3, 12).Value = 1.53
$ws.Cells.Item(3, 13).Value = 2.5
$ws.Cells.Item(3, 14).Value = 2.7
$ws.Cells.Item(3, 15).Value = 1.41
$ws.Cells.Item(3, 16).Value = 1.57
$ws.Cells.Item(3, 17).Value = 2.25
$ws.Cells.Item(3, 18).Value = 2.1
$ws.Cells.Item(3, 19).Value = 1.67
$ws.Cells.Item(3, 20).Value = 6.5
$ws.Cells.Item(3, 22).Value = 11
$ws.Cells.Item(3, 23).Value = 26
$ws.Cells.Item(3, 24).Value = 26
$ws.Cells.Item(3, 25).Value = 41
$ws.Cells.Item(3, 26).Value = 6.5
$ws.Cells.Item(3, 28).Value = 19
$ws.Cells.Item(3, 30).Value = 501
$ws.Cells.Item(3, 31).Value = 7
$ws.Cells.Item(3, 33).Value = 12
$ws.Cells.Item(3, 34).Value = 34
$ws.Cells.Item(3, 35).Value = 29
$ws.Cells.Item(4, 12).Value = 1.36
$ws.Cells.Item(4, 13).Value = 3.2
$ws.Cells.Item(4, 14).Value = 2.05
$ws.Cells.Item(4, 15).Value = 1.68
$ws.Cells.Item(4, 16).Value = 1.44
$ws.Cells.Item(4, 17).Value = 2.63
$ws.Cells.Item(4, 20).Value = 6
$ws.Cells.Item(4, 28).Value = 21
$ws.Cells.Item(4, 30).Value = 451
$ws.Cells.Item(5, 8).Value = 4.1
$ws.Cells.Item(5, 10).Value = 1.05
$ws.Cells.Item(5, 11).Value = 11
$ws.Cells.Item(5, 12).Value = 1.29
$ws.Cells.Item(5, 13).Value = 3.75
$ws.Cells.Item(5, 14).Value = 1.85
$ws.Cells.Item(5, 16).Value = 1.36
$ws.Cells.Item(5, 17).Value = 3
$ws.Cells.Item(5, 18).Value = 1.83
$ws.Cells.Item(5, 19).Value = 1.83
$ws.Cells.Item(5, 20).Value = 7
$ws.Cells.Item(5, 21).Value = 7.5
$ws.Cells.Item(5, 25).Value = 26
$ws.Cells.Item(5, 26).Value = 12
$ws.Cells.Item(5, 27).Value = 7.5
$ws.Cells.Item(5, 31).Value = 13
$ws.Cells.Item(5, 33).Value = 17
$ws.Cells.Item(6, 7).Value = 1.36
$ws.Cells.Item(6, 8).Value = 4.6
$ws.Cells.Item(6, 9).Value = 6.4
$ws.Cells.Item(6, 20).Value = 8.5
$ws.Cells.Item(6, 21).Value = 7.2
$ws.Cells.Item(6, 22).Value = 7.2
$ws.Cells.Item(6, 24).Value = 8.5
$ws.Cells.Item(6, 25).Value = 15.5
$ws.Cells.Item(6, 26).Value = 18.5
$ws.Cells.Item(6, 27).Value = 8.75
$ws.Cells.Item(6, 29).Value = 40
$ws.Cells.Item(6, 30).Value = 200
$ws.Cells.Item(6, 31).Value = 21
$ws.Cells.Item(6, 32).Value = 40
$ws.Cells.Item(6, 33).Value = 17
$ws.Cells.Item(6, 34).Value = 110
$ws.Cells.Item(6, 36).Value = 35
$ws.Cells.Item(7, 18).Value = 1.78
$ws.Cells.Item(7, 19).Value = 1.92
$ws.Cells.Item(8, 9).Value = 3.25
$ws.Cells.Item(8, 12).Value = 1.2
$ws.Cells.Item(8, 13).Value = 4.2
$ws.Cells.Item(8, 18).Value = 1.52
$ws.Cells.Item(8, 19).Value = 2.35
$ws.Cells.Item(8, 21).Value = 11
$ws.Cells.Item(8, 29).Value = 41
$ws.Cells.Item(8, 34).Value = 41
$ws.Cells.Item(8, 36).Value = 29
$ws.Cells.Item(9, 9).Value = 3.7
$ws.Cells.Item(9, 16).Value = 1.38
$ws.Cells.Item(10, 14).Value = 1.65
$ws.Cells.Item(10, 15).Value = 2.2
$ws.Cells.Item(11, 7).Value = 1.73
$ws.Cells.Item(11, 9).Value = 4.1
$ws.Cells.Item(11, 10).Value = 1.05
$ws.Cells.Item(11, 11).Value = 8.75
$ws.Cells.Item(11, 12).Value = 1.24
$ws.Cells.Item(11, 13).Value = 3.75
$ws.Cells.Item(11, 18).Value = 1.7
$ws.Cells.Item(11, 19).Value = 2.05
$ws.Cells.Item(11, 21).Value = 9
$ws.Cells.Item(12, 7).Value = 2.35
$ws.Cells.Item(12, 8).Value = 3.1
$ws.Cells.Item(12, 11).Value = 7.4
$ws.Cells.Item(12, 12).Value = 1.32
$ws.Cells.Item(12, 13).Value = 3.2
$ws.Cells.Item(12, 14).Value = 2.05
$ws.Cells.Item(12, 15).Value = 1.8
$ws.Cells.Item(12, 16).Value = 1.39
$ws.Cells.Item(12, 17).Value = 2.85
$ws.Cells.Item(12, 19).Value = 2.02
$ws.Cells.Item(12, 21).Value = 12
$ws.Cells.Item(12, 24).Value = 21
$ws.Cells.Item(12, 27).Value = 6
$ws.Cells.Item(12, 28).Value = 13
$ws.Cells.Item(12, 30).Value = 251
$ws.Cells.Item(13, 3).Value = "13:00"
$ws.Cells.Item(13, 8).Value = 3.15
$ws.Cells.Item(13, 9).Value = 3
$ws.Cells.Item(13, 11).Value = 6.4
$ws.Cells.Item(13, 18).Value = 1.9
$ws.Cells.Item(13, 19).Value = 1.82
$ws.Cells.Item(13, 20).Value = 6.8
$ws.Cells.Item(13, 21).Value = 10
$ws.Cells.Item(13, 22).Value = 9.25
$ws.Cells.Item(13, 25).Value = 35
$ws.Cells.Item(13, 26).Value = 6.4
$ws.Cells.Item(13, 27).Value = 6.2
$ws.Cells.Item(13, 28).Value = 16
$ws.Cells.Item(13, 30).Value = 800
$ws.Cells.Item(13, 32).Value = 15
$ws.Cells.Item(13, 35).Value = 29
$ws.Cells.Item(14, 7).Value = 1.6
$ws.Cells.Item(14, 9).Value = 5.25
$ws.Cells.Item(14, 10).Value = 1.04
$ws.Cells.Item(14, 11).Value = 13
$ws.Cells.Item(14, 18).Value = 1.75
$ws.Cells.Item(14, 19).Value = 2
$ws.Cells.Item(14, 20).Value = 8
$ws.Cells.Item(14, 21).Value = 8
$ws.Cells.Item(14, 23).Value = 12
$ws.Cells.Item(14, 25).Value = 23
$ws.Cells.Item(14, 26).Value = 13
$ws.Cells.Item(14, 30).Value = 201
$ws.Cells.Item(14, 31).Value = 17
$ws.Cells.Item(14, 33).Value = 17
$ws.Cells.Item(16, 8).Value = 7.1
$ws.Cells.Item(16, 9).Value = 26
$ws.Cells.Item(16, 14).Value = 1.26
$ws.Cells.Item(16, 18).Value = 2.24
$ws.Cells.Item(16, 19).Value = 1.58
$ws.Cells.Item(16, 21).Value = 6.4
$ws.Cells.Item(16, 22).Value = 10.5
$ws.Cells.Item(16, 23).Value = 5.7
$ws.Cells.Item(16, 24).Value = 9.25
$ws.Cells.Item(16, 25).Value = 27
$ws.Cells.Item(16, 26).Value = 21
$ws.Cells.Item(16, 27).Value = 17
$ws.Cells.Item(16, 28).Value = 32
$ws.Cells.Item(16, 35).Value = 500
$ws.Cells.Item(17, 16).Value = 1.44
$ws.Cells.Item(17, 17).Value = 2.63
$ws.Cells.Item(18, 7).Value = 1.3
$ws.Cells.Item(18, 8).Value = 4.75
$ws.Cells.Item(18, 9).Value = 8.5
$ws.Cells.Item(18, 18).Value = 1.8
$ws.Cells.Item(18, 19).Value = 1.91
$ws.Cells.Item(18, 23).Value = 9
$ws.Cells.Item(18, 32).Value = 41
$ws.Cells.Item(19, 7).Value = 1.7
$ws.Cells.Item(19, 8).Value = 3.7
$ws.Cells.Item(19, 9).Value = 4.2
$ws.Cells.Item(19, 10).Value = 1.04
$ws.Cells.Item(19, 11).Value = 8.25
$ws.Cells.Item(19, 12).Value = 1.23
$ws.Cells.Item(19, 13).Value = 3.75
$ws.Cells.Item(19, 14).Value = 1.7
$ws.Cells.Item(19, 15).Value = 2.1
$ws.Cells.Item(19, 16).Value = 1.35
$ws.Cells.Item(19, 17).Value = 2.95
$ws.Cells.Item(19, 21).Value = 9
$ws.Cells.Item(19, 24).Value = 13
$ws.Cells.Item(19, 27).Value = 7.5
$ws.Cells.Item(19, 28).Value = 15
$ws.Cells.Item(19, 31).Value = 15
$ws.Cells.Item(19, 32).Value = 23
$ws.Cells.Item(19, 33).Value = 15
$ws.Cells.Item(19, 34).Value = 51
$ws.Cells.Item(19, 35).Value = 34
$ws.Cells.Item(20, 7).Value = 3.7
$ws.Cells.Item(20, 8).Value = 3.25
$ws.Cells.Item(20, 9).Value = 1.91
$ws.Cells.Item(20, 20).Value = 11
$ws.Cells.Item(20, 21).Value = 19
$ws.Cells.Item(20, 22).Value = 13
$ws.Cells.Item(20, 23).Value = 41
$ws.Cells.Item(20, 25).Value = 41
$ws.Cells.Item(20, 26).Value = 10
$ws.Cells.Item(20, 27).Value = 6.5
$ws.Cells.Item(20, 28).Value = 15
$ws.Cells.Item(20, 30).Value = 251
$ws.Cells.Item(20, 31).Value = 7.5
$ws.Cells.Item(20, 32).Value = 9.5
$ws.Cells.Item(20, 34).Value = 17
$ws.Cells.Item(20, 36).Value = 26
$ws.Cells.Item(22, 7).Value = 1.57
$ws.Cells.Item(22, 8).Value = 3.75
$ws.Cells.Item(22, 9).Value = 5.6
$ws.Cells.Item(22, 11).Value = 7.3
$ws.Cells.Item(22, 12).Value = 1.3
$ws.Cells.Item(22, 13).Value = 3.2
$ws.Cells.Item(22, 14).Value = 1.88
$ws.Cells.Item(22, 15).Value = 1.82
$ws.Cells.Item(22, 17).Value = 2.65
$ws.Cells.Item(22, 18).Value = 1.93
$ws.Cells.Item(22, 19).Value = 1.78
$ws.Cells.Item(22, 21).Value = 7
$ws.Cells.Item(22, 23).Value = 11.25
$ws.Cells.Item(22, 25).Value = 28
$ws.Cells.Item(22, 26).Value = 7.3
$ws.Cells.Item(22, 27).Value = 7.3
$ws.Cells.Item(22, 28).Value = 17.5
$ws.Cells.Item(22, 30).Value = 700
$ws.Cells.Item(22, 33).Value = 17.5
$ws.Cells.Item(22, 36).Value = 60
$ws.Cells.Item(23, 7).Value = 2.2
$ws.Cells.Item(23, 8).Value = 3.25
$ws.Cells.Item(23, 9).Value = 3.4
$ws.Cells.Item(23, 21).Value = 10
$ws.Cells.Item(23, 22).Value = 9.5
$ws.Cells.Item(23, 23).Value = 21
$ws.Cells.Item(23, 24).Value = 19
$ws.Cells.Item(23, 27).Value = 6
$ws.Cells.Item(23, 33).Value = 12
$ws.Cells.Item(23, 34).Value = 34
